$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale row-outline grouping metadata (outlineLevelRow) left over
# from earlier editing; no rows currently carry an outline level. Scope the
# call to the existing used range so it doesn't force-materialize the whole
# 1,048,576-row grid.
$ws.Range("A1:A4").ClearOutline()

# Replace the old placeholder id value and append the new PlayerCardBag ids.
$ws.Range("A4").Value = "10001"
$ws.Range("A5").Value = "10002"
$ws.Range("A6").Value = "10003"
$ws.Range("A7").Value = "10004"
$ws.Range("A8").Value = "10005"
$ws.Range("A9").Value = "10006"
$ws.Range("A10").Value = "10007"

# Match the author's final cursor position.
$ws.Range("B14").Select()
